$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.3912683333333333
$ws.Range("H2").Value = 1.173805
$ws.Range("I2").Value = 0.004923718964983145
$ws.Range("J2").Value = 0.004923718964983145
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.779790999999999
$ws.Range("N2").Value = 14.339373
$ws.Range("Q2").Value = 1.870180858251666
$ws.Range("R2").Value = 16.831627724265
$ws.Range("S2").Value = 0.004923718964983145
$ws.Range("T2").Value = 0.004923718964983145

# Row 3
$ws.Range("I3").Value = 0.07888477275715973
$ws.Range("J3").Value = 0.07888477275715973
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.779790999999999
$ws.Range("N3").Value = 14.339373
$ws.Range("Q3").Value = 29.96287827700533
$ws.Range("R3").Value = 269.665904493048
$ws.Range("S3").Value = 0.07888477275715973
$ws.Range("T3").Value = 0.07888477275715973

# Row 4
$ws.Range("G4").Value = 72.68848166666666
$ws.Range("H4").Value = 218.065445
$ws.Range("I4").Value = 0.9147115297293749
$ws.Range("J4").Value = 0.9147115297293749
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.779790999999999
$ws.Range("N4").Value = 14.339373
$ws.Range("Q4").Value = 347.4357504739983
$ws.Range("R4").Value = 3126.921754265984
$ws.Range("S4").Value = 0.9147115297293749
$ws.Range("T4").Value = 0.9147115297293749

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.117608
$ws.Range("H5").Value = 0.352824
$ws.Range("I5").Value = 0.001479978548482255
$ws.Range("J5").Value = 0.001479978548482255
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.779790999999999
$ws.Range("N5").Value = 14.339373
$ws.Range("Q5").Value = 0.5621416599279999
$ws.Range("R5").Value = 5.059274939352
$ws.Range("S5").Value = 0.001479978548482255
$ws.Range("T5").Value = 0.001479978548482255
